$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AA6" = 0.8172043010752689
    "AB6" = 0.7956989247311828
    "AD6" = 0.7940288263555251
    "AE6" = 0.02357546057468482
    "AF6" = 0.776595744680851
    "AG6" = 0.7849462365591398
    "AI6" = 0.7741935483870968
    "AJ6" = 0.8387096774193549
    "AK6" = 0.8219629375428964
    "AL6" = 0.04851581611932004
    "AN6" = 0.7526881720430108
    "AP6" = 0.8279569892473119
    "B6" = 0.8305879661404714
    "C6" = 0.0361509544709699
    "D6" = 0.776595744680851
    "E6" = 0.8064516129032258
    "F6" = 0.8494623655913979
    "G6" = 0.8387096774193549
    "H6" = 0.8817204301075269
    "I6" = 0.8112331274307939
    "J6" = 0.03653407699835858
    "K6" = 0.776595744680851
    "L6" = 0.7634408602150538
    "M6" = 0.8387096774193549
    "O6" = 0.8602150537634409
    "P6" = 0.7253717684740335
    "Q6" = 0.02626879524596896
    "R6" = 0.7021276595744681
    "S6" = 0.7096774193548387
    "T6" = 0.7311827956989247
    "U6" = 0.7096774193548387
    "V6" = 0.7741935483870968
    "W6" = 0.7962251201098148
    "X6" = 0.02548656950758862
    "Y6" = 0.7553191489361702
    "Z6" = 0.7849462365591398
    "AC7" = 0.8709677419354839
    "AD7" = 0.8562113932738503
    "AE7" = 0.03511569808548161
    "AF7" = 0.8617021276595744
    "AH7" = 0.8709677419354839
    "AJ7" = 0.9032258064516129
    "AK7" = 0.8369251887439946
    "AL7" = 0.02477828029829682
    "AM7" = 0.8297872340425532
    "AN7" = 0.8064516129032258
    "B7" = 0.8433539235872798
    "C7" = 0.04976362478805044
    "F7" = 0.8709677419354839
    "I7" = 0.8584305650880806
    "J7" = 0.03111597825317717
    "M7" = 0.8602150537634409
    "N7" = 0.8602150537634409
    "P7" = 0.8283687943262411
    "Q7" = 0.04643950885054888
    "R7" = 0.8085106382978723
    "S7" = 0.7849462365591398
    "T7" = 0.8387096774193549
    "V7" = 0.9139784946236559
    "W7" = 0.8455044612216884
    "X7" = 0.03367115730482832
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
